$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix inconsistent "o"/progress-text values in column D to uppercase "O"
$ws.Range("D4").Value = "O"
$ws.Range("D5").Value = "O"
$ws.Range("D6").Value = "O"
$ws.Range("D7").Value = "O"

# Add lowercase "o" markers to column D for rows 8-10
$ws.Range("D8").Value = "o"
$ws.Range("D9").Value = "o"
$ws.Range("D10").Value = "o"

# Add "O" markers to new column H for rows 2,3,4,5,6
$ws.Range("H2").Value = "O"
$ws.Range("H3").Value = "O"
$ws.Range("H4").Value = "O"
$ws.Range("H5").Value = "O"
$ws.Range("H6").Value = "O"
